$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 36; this shifts rows 36..85 down to 37..86
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the new data record
$ws.Cells.Item(36, 1).Value = 10
$ws.Cells.Item(36, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value = "La Araucanía"
$ws.Cells.Item(36, 4).Value2 = 45028
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100107
$ws.Cells.Item(36, 8).Value = "Otros"
$ws.Cells.Item(36, 9).Value = 100107011
$ws.Cells.Item(36, 10).Value = "Tuna"
$ws.Cells.Item(36, 11).Value = "Sin especificar"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 60
$ws.Cells.Item(36, 14).Value = 22000
$ws.Cells.Item(36, 15).Value = 22000
$ws.Cells.Item(36, 16).Value = 22000
$ws.Cells.Item(36, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(36, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(36, 19).Value = 1375
$ws.Cells.Item(36, 20).Value = 16
